$d = $word.ActiveDocument

# --- 1. Insert the new "Meta description" paragraph right after the H1 title ---
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$metaEnd = $metaPara.Range.End
$metaFullText = "Meta description: Experience epic wins and 4 progressive jackpots. Read our comprehensive review of the Age of the Gods: Goddess of Wisdom slot, and play for free."
$metaFullRange = $d.Range($metaStart, $metaEnd - 1)
$metaFullRange.Text = $metaFullText

# Bold just the "Meta description" label (leave the rest, incl. the colon, un-bold)
$labelRange = $d.Range($metaStart, $metaStart + 16)
$labelRange.Bold = 1

# --- 2. Remove the duplicated bold title paragraph near the end of the document ---
$duplicateParaFound = $false
For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    if ($i -gt 2 -and $txt -eq "Play Age of the Gods: Goddess of Wisdom Slot for Free | Review") {
        $para.Range.Delete()
        $duplicateParaFound = $true
        break
    }
}

# --- 3. Replace the italic meta-description paragraph text with the image prompt ---
$promptText = 'Create a feature image for "Age of the Gods: Goddess of Wisdom". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be standing in front of an ancient Greek temple while holding a shield with the game title written on it. Athena, the central figure of the game, should be standing next to the warrior with a confident stance. The symbols of the game, including the Gorgoneion, Olive Branches, Helmets, and Armor, should be seen floating around the two figures. The image should be colorful and eye-catching to attract players'' attention.'

For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $txt = $para.Range.Text.TrimEnd([char]13)
    if ($txt -eq "Experience epic wins and 4 progressive jackpots. Read our comprehensive review of the Age of the Gods: Goddess of Wisdom slot, and play for free.") {
        $pStart = $para.Range.Start
        $pEnd = $para.Range.End
        $fullRange = $d.Range($pStart, $pEnd - 1)
        $fullRange.Text = $promptText
        break
    }
}

Write-Output "done"
